$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: insert new "Jan_2026" column, roll Dec_2025/Nov_2025 forward, drop Oct_2025
$ws.Range("D1").Value = "Jan_2026"
$ws.Range("E1").Value = "Dec_2025"
$ws.Range("F1").Value = "Nov_2025"

# Row 2: Piramal Finance Ltd
$ws.Range("A2").Value = "INE202B01038"
$ws.Range("B2").Value = "Piramal Finance Ltd"
$ws.Range("C2").Value = "quant BFSI Fund"
$ws.Range("D2").Value = 9.314091
$ws.Range("E2").Value = 8.685552
$ws.Range("F2").Value = 8.633141
$ws.Range("G2").Value = 0.628539
$ws.Range("H2").Value = 0.6809499999999993

# Row 3: Shriram Finance Limited
$ws.Range("A3").Value = "INE721A01047"
$ws.Range("B3").Value = "Shriram Finance Limited"
$ws.Range("C3").Value = "quant BFSI Fund"
$ws.Range("D3").Value = 8.767368
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 8.767368
$ws.Range("H3").Value = 8.767368

# Row 4: Capri Global Capital Limited
$ws.Range("A4").Value = "INE180C01042"
$ws.Range("B4").Value = "Capri Global Capital Limited"
$ws.Range("C4").Value = "quant BFSI Fund"
$ws.Range("D4").Value = 8.001052
$ws.Range("E4").Value = 7.774914
$ws.Range("F4").Value = 8.436206
$ws.Range("G4").Value = 0.2261379999999997
$ws.Range("H4").Value = -0.4351540000000007

# Row 5: HDFC Bank Limited
$ws.Range("A5").Value = "INE040A01034"
$ws.Range("B5").Value = "HDFC Bank Limited"
$ws.Range("C5").Value = "quant BFSI Fund"
$ws.Range("D5").Value = 7.336848
$ws.Range("E5").Value = 0.471945
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 6.864903
$ws.Range("H5").Value = 7.336848

# Row 6: LIC Housing Finance Ltd
$ws.Range("A6").Value = "INE115A01026"
$ws.Range("B6").Value = "LIC Housing Finance Ltd"
$ws.Range("C6").Value = "quant BFSI Fund"
$ws.Range("D6").Value = 6.249933
$ws.Range("E6").Value = 6.340376
$ws.Range("F6").Value = 6.80419
$ws.Range("G6").Value = -0.09044299999999961
$ws.Range("H6").Value = -0.5542569999999998

# Row 7: Kotak Mahindra Bank Limited
$ws.Range("A7").Value = "INE237A01036"
$ws.Range("B7").Value = "Kotak Mahindra Bank Limited"
$ws.Range("C7").Value = "quant BFSI Fund"
$ws.Range("D7").Value = 5.747081
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 5.747081
$ws.Range("H7").Value = 5.747081

# Row 8: Bajaj Finance Limited
$ws.Range("A8").Value = "INE296A01032"
$ws.Range("B8").Value = "Bajaj Finance Limited"
$ws.Range("C8").Value = "quant BFSI Fund"
$ws.Range("D8").Value = 5.264682
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 3.416392
$ws.Range("G8").Value = 5.264682
$ws.Range("H8").Value = 1.84829

# Row 9: ICICI Prudential AMC Ltd
$ws.Range("A9").Value = "INE346A01027"
$ws.Range("B9").Value = "ICICI Prudential AMC Ltd"
$ws.Range("C9").Value = "quant BFSI Fund"
$ws.Range("D9").Value = 5.021502
$ws.Range("E9").Value = 4.448358
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0.5731440000000001
$ws.Range("H9").Value = 5.021502

# Row 10: HDFC Asset Management Company Ltd
$ws.Range("A10").Value = "INE127D01025"
$ws.Range("B10").Value = "HDFC Asset Management Company Ltd"
$ws.Range("C10").Value = "quant BFSI Fund"
$ws.Range("D10").Value = 4.974659
$ws.Range("E10").Value = 3.728754
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 1.245905
$ws.Range("H10").Value = 4.974659

# Row 11: HDFC Life Insurance Co Ltd
$ws.Range("A11").Value = "INE795G01014"
$ws.Range("B11").Value = "HDFC Life Insurance Co Ltd"
$ws.Range("C11").Value = "quant BFSI Fund"
$ws.Range("D11").Value = 4.448754
$ws.Range("E11").Value = 3.723976
$ws.Range("F11").Value = 2.840253
$ws.Range("G11").Value = 0.7247780000000001
$ws.Range("H11").Value = 1.608501

# Row 12: Adani Enterprises Limited
$ws.Range("A12").Value = "INE423A01024"
$ws.Range("B12").Value = "Adani Enterprises Limited"
$ws.Range("C12").Value = "quant BFSI Fund"
$ws.Range("D12").Value = 3.418112
$ws.Range("E12").Value = 3.749559
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = -0.3314470000000003
$ws.Range("H12").Value = 3.418112

# Row 13: Nippon Life India Asset Management Ltd
$ws.Range("A13").Value = "INE298J01013"
$ws.Range("B13").Value = "Nippon Life India Asset Management Ltd"
$ws.Range("C13").Value = "quant BFSI Fund"
$ws.Range("D13").Value = 2.141619
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 2.141619
$ws.Range("H13").Value = 2.141619

# Row 14: ICICI Bank Limited
$ws.Range("A14").Value = "INE090A01021"
$ws.Range("B14").Value = "ICICI Bank Limited"
$ws.Range("C14").Value = "quant BFSI Fund"
$ws.Range("D14").Value = 0.489655
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0.489655
$ws.Range("H14").Value = 0.489655

# Row 15: ICICI Prudential Life Insurance Co Ltd
$ws.Range("A15").Value = "INE726G01019"
$ws.Range("B15").Value = "ICICI Prudential Life Insurance Co Ltd"
$ws.Range("C15").Value = "quant BFSI Fund"
$ws.Range("D15").Value = 0.039149
$ws.Range("E15").Value = 2.836314
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = -2.797165
$ws.Range("H15").Value = 0.039149

# Row 16: Reliance Industries Limited
$ws.Range("A16").Value = "INE002A01018"
$ws.Range("B16").Value = "Reliance Industries Limited"
$ws.Range("C16").Value = "quant BFSI Fund"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1.551008
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = -1.551008

# Row 17: Anand Rathi Wealth Limited
$ws.Range("A17").Value = "INE463V01026"
$ws.Range("B17").Value = "Anand Rathi Wealth Limited"
$ws.Range("C17").Value = "quant BFSI Fund"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 7.22842
$ws.Range("F17").Value = 7.074929
$ws.Range("G17").Value = -7.22842
$ws.Range("H17").Value = -7.074929

# Row 18: Kotak Mahindra Bank Limited
$ws.Range("A18").Value = "INE237A01028"
$ws.Range("B18").Value = "Kotak Mahindra Bank Limited"
$ws.Range("C18").Value = "quant BFSI Fund"
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 3.139423
$ws.Range("F18").Value = 3.195126
$ws.Range("G18").Value = -3.139423
$ws.Range("H18").Value = -3.195126

# Row 19: SBI Cards & Payment Services Ltd
$ws.Range("A19").Value = "INE018E01016"
$ws.Range("B19").Value = "SBI Cards & Payment Services Ltd"
$ws.Range("C19").Value = "quant BFSI Fund"
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 6.086461
$ws.Range("F19").Value = 5.106652
$ws.Range("G19").Value = -6.086461
$ws.Range("H19").Value = -5.106652

# Row 20: Life Insurance Corporation Of India
$ws.Range("A20").Value = "INE0J1Y01017"
$ws.Range("B20").Value = "Life Insurance Corporation Of India"
$ws.Range("C20").Value = "quant BFSI Fund"
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 8.064005
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = -8.064005

# Row 21: Canara HSBC Life Insurance Company Ltd
$ws.Range("A21").Value = "INE01TY01017"
$ws.Range("B21").Value = "Canara HSBC Life Insurance Company Ltd"
$ws.Range("C21").Value = "quant BFSI Fund"
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 10.298525
$ws.Range("F21").Value = 9.137125
$ws.Range("G21").Value = -10.298525
$ws.Range("H21").Value = -9.137125

# Row 22: Bajaj Finserv Ltd.
$ws.Range("A22").Value = "INE918I01026"
$ws.Range("B22").Value = "Bajaj Finserv Ltd."
$ws.Range("C22").Value = "quant BFSI Fund"
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 1.978361
$ws.Range("F22").Value = 2.141484
$ws.Range("G22").Value = -1.978361
$ws.Range("H22").Value = -2.141484
